$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1179
$ws.Range("J32").Value = 975.5
$ws.Range("L32").Value = 975.5
$ws.Range("N32").Value = -1627.5
$ws.Range("H62").Value = 47746.5
$ws.Range("I62").Value = 17249.572
$ws.Range("K62").Value = 17249.572
$ws.Range("M62").Value = -16625.572
$ws.Range("H65").Value = 47746.5
$ws.Range("I65").Value = 17249.572
$ws.Range("K65").Value = 86247.86
$ws.Range("M65").Value = -83127.86
$ws.Range("H87").Value = 49998
$ws.Range("J87").Value = 49998
$ws.Range("L87").Value = 49998
$ws.Range("N87").Value = -52494
$ws.Range("H90").Value = 49998
$ws.Range("J90").Value = 49998
$ws.Range("L90").Value = 149994
$ws.Range("N90").Value = -162474
$ws.Range("H132").Value = 1784.289
$ws.Range("I132").Value = 1767.6904
$ws.Range("K132").Value = 5303.0712
$ws.Range("M132").Value = -2773.0712

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4262970
$ws.Range("I32").Value = 4883296
$ws.Range("K32").Value = 4883296
$ws.Range("M32").Value = -4883009
$ws.Range("H45").Value = 2952.3
$ws.Range("I45").Value = 2912
$ws.Range("K45").Value = 2912
$ws.Range("M45").Value = -2535
$ws.Range("H61").Value = 4237.675
$ws.Range("I61").Value = 1469.7576
$ws.Range("K61").Value = 1469.7576
$ws.Range("M61").Value = -1257.7576
$ws.Range("H74").Value = 2477.6177
$ws.Range("I74").Value = 1627.9231
$ws.Range("K74").Value = 1627.9231
$ws.Range("M74").Value = -753.9231
$ws.Range("H77").Value = 2477.6177
$ws.Range("I77").Value = 1627.9231
$ws.Range("K77").Value = 8139.6155
$ws.Range("M77").Value = -3771.6155
$ws.Range("H82").Value = 47181
$ws.Range("J82").Value = 47181
$ws.Range("L82").Value = 47181
$ws.Range("N82").Value = -47903
$ws.Range("H85").Value = 47181
$ws.Range("J85").Value = 47181
$ws.Range("L85").Value = 47181
$ws.Range("N85").Value = -49677
$ws.Range("H97").Value = 2604840.5
$ws.Range("I97").Value = 518.3333
$ws.Range("J97").Value = 41669670
$ws.Range("K97").Value = 518.3333
$ws.Range("L97").Value = 41669670
$ws.Range("M97").Value = -22.33330000000001
$ws.Range("N97").Value = -41670662
$ws.Range("H132").Value = 5481.237
$ws.Range("I132").Value = 2309
$ws.Range("K132").Value = 6927
$ws.Range("M132").Value = -4397
$ws.Range("H136").Value = 4237.675
$ws.Range("I136").Value = 1469.7576
$ws.Range("K136").Value = 4409.2728
$ws.Range("M136").Value = -1859.2728

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 113730.22
$ws.Range("I86").Value = 201199
$ws.Range("K86").Value = 201199
$ws.Range("M86").Value = -200076
$ws.Range("H89").Value = 113730.22
$ws.Range("I89").Value = 201199
$ws.Range("K89").Value = 1005995
$ws.Range("M89").Value = -1000379
$ws.Range("H94").Value = 1555.3846
$ws.Range("I94").Value = 974.7727
$ws.Range("J94").Value = 4748.75
$ws.Range("K94").Value = 974.7727
$ws.Range("L94").Value = 4748.75
$ws.Range("M94").Value = -523.7727
$ws.Range("N94").Value = -5650.75
$ws.Range("H134").Value = 6955.517
$ws.Range("I134").Value = 2600.6667
$ws.Range("J134").Value = 11621.429
$ws.Range("K134").Value = 7802.000100000001
$ws.Range("L134").Value = 34864.287
$ws.Range("M134").Value = -5267.000100000001
$ws.Range("N134").Value = -39934.287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 83333440
$ws.Range("I7").Value = 75.5
$ws.Range("J7").Value = 250000180
$ws.Range("K7").Value = 75.5
$ws.Range("L7").Value = 250000180
$ws.Range("M7").Value = 37.5
$ws.Range("N7").Value = -250000406
$ws.Range("H13").Value = 502
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("N13").Value = -1278
$ws.Range("H16").Value = 3830.6086
$ws.Range("J16").Value = 6267.6665
$ws.Range("L16").Value = 6267.6665
$ws.Range("N16").Value = -6841.6665
$ws.Range("H94").Value = 1321.6154
$ws.Range("J94").Value = 1049.8889
$ws.Range("L94").Value = 1049.8889
$ws.Range("N94").Value = -1951.8889
$ws.Range("H105").Value = 14286357
$ws.Range("I105").Value = 17857744
$ws.Range("K105").Value = 17857744
$ws.Range("M105").Value = -17855997
$ws.Range("H113").Value = 3830.6086
$ws.Range("J113").Value = 6267.6665
$ws.Range("L113").Value = 6267.6665
$ws.Range("N113").Value = -10607.6665
$ws.Range("H132").Value = 5678.6772
$ws.Range("I132").Value = 3051.8667
$ws.Range("J132").Value = 8141.3125
$ws.Range("K132").Value = 9155.6001
$ws.Range("L132").Value = 24423.9375
$ws.Range("M132").Value = -6625.6001
$ws.Range("N132").Value = -29483.9375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 131890.8
$ws.Range("I2").Value = 25044.084
$ws.Range("K2").Value = 150264.504
$ws.Range("M2").Value = -150151.504
$ws.Range("H46").Value = 250001000
$ws.Range("J46").Value = 500000000
$ws.Range("L46").Value = 1500000000
$ws.Range("N46").Value = -1500000182
$ws.Range("H122").Value = 1572080.6
$ws.Range("J122").Value = 647.9167
$ws.Range("L122").Value = 5831.2503
$ws.Range("N122").Value = -10731.2503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3459.9375
$ws.Range("J80").Value = 5995.4
$ws.Range("L80").Value = 5995.4
$ws.Range("N80").Value = -7991.4
$ws.Range("H83").Value = 3459.9375
$ws.Range("J83").Value = 5995.4
$ws.Range("L83").Value = 29977
$ws.Range("N83").Value = -39961
$ws.Range("H132").Value = 5693.1333
$ws.Range("I132").Value = 1399.7273
$ws.Range("K132").Value = 4199.1819
$ws.Range("M132").Value = -1669.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 37666.668
$ws.Range("I56").Value = 37666.668
$ws.Range("K56").Value = 37666.668
$ws.Range("M56").Value = -36975.668
$ws.Range("H82").Value = 1453.8
$ws.Range("I82").Value = 672.5
$ws.Range("K82").Value = 672.5
$ws.Range("M82").Value = -311.5
$ws.Range("H85").Value = 1453.8
$ws.Range("I85").Value = 672.5
$ws.Range("K85").Value = 672.5
$ws.Range("M85").Value = 575.5
$ws.Range("H93").Value = 982.4286
$ws.Range("I93").Value = 895.5
$ws.Range("J93").Value = 1199.75
$ws.Range("K93").Value = 895.5
$ws.Range("L93").Value = 1199.75
$ws.Range("M93").Value = 352.5
$ws.Range("N93").Value = -3695.75
$ws.Range("H94").Value = 43640
$ws.Range("J94").Value = 43640
$ws.Range("L94").Value = 43640
$ws.Range("N94").Value = -44992
$ws.Range("H97").Value = 33784.332
$ws.Range("J97").Value = 33784.332
$ws.Range("L97").Value = 33784.332
$ws.Range("N97").Value = -35766.332
$ws.Range("H132").Value = 10195.268
$ws.Range("I132").Value = 5382.885
$ws.Range("K132").Value = 16148.655
$ws.Range("M132").Value = -13618.655
$ws.Range("H136").Value = 12900.714
$ws.Range("I136").Value = 3385.6843
$ws.Range("J136").Value = 20760.957
$ws.Range("K136").Value = 10157.0529
$ws.Range("L136").Value = 62282.871
$ws.Range("M136").Value = -7607.052899999999
$ws.Range("N136").Value = -67382.871

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 3073
$ws.Range("H74").Value = 6250
$ws.Range("J74").Value = 6250
$ws.Range("L74").Value = 6250
$ws.Range("N74").Value = -8122
$ws.Range("H77").Value = 6250
$ws.Range("J77").Value = 6250
$ws.Range("L77").Value = 18750
$ws.Range("N77").Value = -28110
$ws.Range("H132").Value = 62509610
$ws.Range("I132").Value = 100010970
$ws.Range("K132").Value = 300032910
$ws.Range("M132").Value = -300030380
$ws.Range("H136").Value = 439850.4
$ws.Range("I136").Value = 1504.5454
$ws.Range("K136").Value = 4513.6362
$ws.Range("M136").Value = -1963.6362
